$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.57"
$ws.Range("E2").Value = "'-0.68%"
$ws.Range("D3").Value = "'43.79"
$ws.Range("E3").Value = "'5.73%"
$ws.Range("D4").Value = "'5.441"
$ws.Range("E4").Value = "'-3.34%"
$ws.Range("D5").Value = "'0.08160"
$ws.Range("E5").Value = "'-2.17%"
$ws.Range("D6").Value = "'8.727"
$ws.Range("E6").Value = "'-0.70%"
$ws.Range("D7").Value = "'1.914"
$ws.Range("E7").Value = "'-3.12%"
$ws.Range("D8").Value = "'4.315"
$ws.Range("E8").Value = "'-3.72%"
$ws.Range("D9").Value = "'2.786"
$ws.Range("E9").Value = "'-4.64%"
$ws.Range("D10").Value = "'0.9426"
$ws.Range("E10").Value = "'1.85%"
$ws.Range("D11").Value = "'0.1180"
$ws.Range("E11").Value = "'-8.72%"
$ws.Range("D12").Value = "'0.1896"
$ws.Range("E12").Value = "'-3.13%"
$ws.Range("D13").Value = "'0.09769"
$ws.Range("E13").Value = "'4.33%"
$ws.Range("D14").Value = "'0.04192"
$ws.Range("E14").Value = "'7.19%"
$ws.Range("D15").Value = "'0.1066"
$ws.Range("E15").Value = "'0.68%"
$ws.Range("D16").Value = "'0.001271"
$ws.Range("E16").Value = "'-2.45%"
$ws.Range("D17").Value = "'0.006064"
$ws.Range("E17").Value = "'-0.45%"
$ws.Range("D18").Value = "'3.547"
$ws.Range("E18").Value = "'3.06%"
$ws.Range("D20").Value = "'8.752"
$ws.Range("E20").Value = "'2.58%"
$ws.Range("D21").Value = "'0.1360"
$ws.Range("E21").Value = "'-0.86%"
$ws.Range("D22").Value = "'0.2501"
$ws.Range("E22").Value = "'3.33%"
$ws.Range("D23").Value = "'0.04370"
$ws.Range("E23").Value = "'-0.90%"
$ws.Range("D24").Value = "'0.001239"
$ws.Range("E24").Value = "'-2.76%"
$ws.Range("D25").Value = "'0.004339"
$ws.Range("E25").Value = "'-1.10%"
$ws.Range("D26").Value = "'0.0001235"
$ws.Range("E26").Value = "'3.00%"
$ws.Range("D27").Value = "'0.0004006"
$ws.Range("E27").Value = "'31.55%"
$ws.Range("D39").Value = "'0.02670"
$ws.Range("E39").Value = "'-5.11%"
$ws.Range("D40").Value = "'0.05659"
$ws.Range("E40").Value = "'2.76%"
$ws.Range("D41").Value = "'0.007869"
$ws.Range("E41").Value = "'-0.96%"
$ws.Range("D42").Value = "'0.009775"
$ws.Range("E42").Value = "'4.90%"
$ws.Range("E43").Value = "'-1.81%"
$ws.Range("D44").Value = "'0.002127"
$ws.Range("E44").Value = "'-0.52%"
$ws.Range("D45").Value = "'0.009608"
$ws.Range("E45").Value = "'-13.38%"
$ws.Range("D46").Value = "'0.00007059"
$ws.Range("E46").Value = "'-0.46%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.38%"
$ws.Range("E48").Value = "'0.79%"
$ws.Range("D49").Value = "'0.002280"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("D50").Value = "'0.00002109"
$ws.Range("E50").Value = "'0.38%"
$ws.Range("D51").Value = "'0.0002009"
$ws.Range("E51").Value = "'0.38%"
